$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD:AF) next to the existing "Salary"/
# "Unnamed: 28" columns, recording each team's season record.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the other header cells (bold font, thin border,
# centered alignment) by copying the format from the neighboring header
# cell "AC1" onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins=74, Losses=88, Ties=0) for every player
# row in the table.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
